$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value  = -7.207000000000001
$ws.Range("C3").Value  = -11.9
$ws.Range("A4").Value  = -21.468
$ws.Range("B4").Value  = 7.013000000000001
$ws.Range("C4").Value  = -12.326
$ws.Range("B5").Value  = 6.324999999999999
$ws.Range("D5").Value  = -8.300999999999998
$ws.Range("A6").Value  = -21.383
$ws.Range("B6").Value  = 6.252000000000001
$ws.Range("A7").Value  = -21.179
$ws.Range("A8").Value  = -21.398
$ws.Range("B8").Value  = 6.153
$ws.Range("C9").Value  = -11.775
$ws.Range("C11").Value = -12.642
$ws.Range("C14").Value = -11.607
$ws.Range("A16").Value = -21.212
$ws.Range("B16").Value = 5.896000000000001
$ws.Range("C18").Value = -12.634
$ws.Range("A20").Value = -21.907
$ws.Range("D20").Value = -8.434000000000001
$ws.Range("A21").Value = -21.14
$ws.Range("B22").Value = 6.386000000000001
$ws.Range("C25").Value = -12.642
